$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 44, shifting existing rows 44:65 down to 45:66.
$ws.Rows("44").Insert()

# Populate the newly inserted row 44 with its data.
$ws.Cells.Item(44, 1).Value = 11
$ws.Cells.Item(44, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(44, 3).Value = "Bíobío"
$ws.Cells.Item(44, 4).Value = Get-Date -Year 2021 -Month 10 -Day 20 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(44, 5).Value = 8
$ws.Cells.Item(44, 6).Value = 100112032
$ws.Cells.Item(44, 7).Value = "Zapallo italiano"
$ws.Cells.Item(44, 8).Value = "Sin especificar"
$ws.Cells.Item(44, 9).Value = "Primera"
$ws.Cells.Item(44, 10).Value = 100
$ws.Cells.Item(44, 11).Value = 8000
$ws.Cells.Item(44, 12).Value = 9000
$ws.Cells.Item(44, 13).Value = 8500
$ws.Cells.Item(44, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(44, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(44, 16).Value = 170
$ws.Cells.Item(44, 17).Value = 50
$ws.Cells.Item(44, 18).Value = "Hortaliza"
